$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Decommissioning date" column (U) -------------------------
# Header cell: same look as the other header cells (copy format from T1,
# the previous last header "Automatic").
$ws.Range("U1").Value = "Decommissioning date"
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)  # xlPasteFormats

# Matching empty data cell under the new header, same look as T2.
$ws.Range("T2").Copy()
$ws.Range("U2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# New column should be as wide as the other date columns (R:T).
$ws.Range("U1").EntireColumn.ColumnWidth = $ws.Range("T1").EntireColumn.ColumnWidth

# --- Drop the unused, blank formatting-only rows 3-10 -----------------------
$ws.Rows("3:10").Delete()

Write-Host "done"
